$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.246613025665283
$ws.Range("B1").Value = 1.138303637504578
$ws.Range("C1").Value = 0.9702181220054626
$ws.Range("D1").Value = 1.021341323852539
$ws.Range("E1").Value = 1.003981351852417
